# #5: cash & deposit done
# Fill out the "存款" (deposit) sheet with a proper header row and a fully
# populated data row, matching the other property sheets' layout:
#   bank | deposit_type | currency | owner | total | property_category |
#   category | date | legislator_name | legislator_id | source_file | index

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# --- Header row (row 1) ---------------------------------------------------
$ws.Cells.Item(1, 2).Value = "bank"
$ws.Cells.Item(1, 3).Value = "deposit_type"
$ws.Cells.Item(1, 4).Value = "currency"
$ws.Cells.Item(1, 5).Value = "owner"
$ws.Cells.Item(1, 6).Value = "total"
$ws.Cells.Item(1, 7).Value = "property_category"
$ws.Cells.Item(1, 8).Value = "category"
$ws.Cells.Item(1, 9).Value = "date"
$ws.Cells.Item(1, 10).Value = "legislator_name"
$ws.Cells.Item(1, 11).Value = "legislator_id"
$ws.Cells.Item(1, 12).Value = "source_file"
$ws.Cells.Item(1, 13).Value = "index"

# New header cells (G1:M1) need the same bold/centered/bordered look as the
# existing header cells (B1:F1) -- copy the formatting from an existing
# header cell instead of rebuilding it property-by-property so we don't
# leave behind a pile of one-off intermediate styles.
$headerRange = $ws.Range($ws.Cells.Item(1, 7), $ws.Cells.Item(1, 13))
$ws.Cells.Item(1, 2).Copy() | Out-Null
$headerRange.PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Data row (row 2) ------------------------------------------------------
$ws.Cells.Item(2, 2).Value = "臺灣銀行."
$ws.Cells.Item(2, 3).Value = "定期存款"
$ws.Cells.Item(2, 4).Value = "新臺幣"
$ws.Cells.Item(2, 5).Value = "蔡煌瑯"
$ws.Cells.Item(2, 6).Value = 1615256
$ws.Cells.Item(2, 7).Value = "deposit"
$ws.Cells.Item(2, 8).Value = "normal"

# "2013-12-17" looks like a date, so a plain .Value assignment would get
# auto-converted to a date serial. Enter it as a string-literal formula
# (never date-sniffed) then paste-special just the value back over itself
# to collapse it to a plain text cell, same as the other date-ish text
# cells elsewhere in this workbook.
$dateCell = $ws.Cells.Item(2, 9)
$dateCell.Formula = '="2013-12-17"'
$dateCell.Copy() | Out-Null
$dateCell.PasteSpecial(-4163) | Out-Null

$ws.Cells.Item(2, 10).Value = "蔡煌瑯"
$ws.Cells.Item(2, 11).Value = 752
$ws.Cells.Item(2, 12).Value = "tmpc9fc1"
$ws.Cells.Item(2, 13).Value = 50
$excel.CutCopyMode = 0
